$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 84 (shifts existing rows 84..188 down to 85..189),
# matching the weekly price-report pattern already present in the sheet
# (same market / category / etc. as the surrounding rows).
$ws.Rows.Item(84).Insert()

$ws.Range("A84").Value = 3
$ws.Range("B84").Value = "Femacal de La Calera"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44483
$ws.Range("E84").Value = 5
$ws.Range("F84").Value = 100112039
$ws.Range("G84").Value = "Ciboulette"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 180
$ws.Range("K84").Value = 1500
$ws.Range("L84").Value = 1500
$ws.Range("M84").Value = 1500
$ws.Range("N84").Value = "`$/docena de atados"
$ws.Range("O84").Value = "Provincia de Quillota"
$ws.Range("P84").Value = 500
$ws.Range("Q84").Value = 3
$ws.Range("R84").Value = "Hortaliza"
